$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 375, shifting existing rows 375-402 down to 376-403.
$ws.Rows(375).Insert()

# Populate the newly inserted row 375 with the new weekly record.
$ws.Range("A375").Value = 5
$ws.Range("B375").Value = "Macroferia Regional de Talca"
$ws.Range("C375").Value = "Maule"
$ws.Range("D375").Value = 45013
$ws.Range("E375").Value = 7
$ws.Range("F375").Value = 100112009
$ws.Range("G375").Value = "Acelga"
$ws.Range("H375").Value = "Sin especificar"
$ws.Range("I375").Value = "Primera"
$ws.Range("J375").Value = 400
$ws.Range("K375").Value = 2500
$ws.Range("L375").Value = 2500
$ws.Range("M375").Value = 2500
$ws.Range("N375").Value = "$/docena de atados (4 kilos)"
$ws.Range("O375").Value = "Región del Maule"
$ws.Range("P375").Value = 625
$ws.Range("Q375").Value = 4
$ws.Range("R375").Value = "Hortaliza"
